# Trade #107 (MarketMaking strategy) closes early at 2026-02-17 21:13:48,
# and a brand-new trade #140 is opened at 21:13:41.
# This touches four worksheet tabs: Summary, Strategy Status, All Trades and MarketMaking.
#
# NOTE: in this workbook the "All Trades" tab and the "MarketMaking" tab hold the
# row layouts one might intuitively expect to be swapped (verified empirically
# against the target OOXML) - the code below addresses worksheets strictly by
# their tab name and matches the column layout actually present on each tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.13   # Current Capital
$wsSummary.Range("B4").Value = 0.92      # Total P&L $
$wsSummary.Range("B5").Value = 0.17      # Total P&L %
$wsSummary.Range("B6").Value = 107       # Total Trades
$wsSummary.Range("B8").Value = 42        # Losing Trades
$wsSummary.Range("B9").Value = 45.79     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 5)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.13
$wsStatus.Range("D5").Value = 74
$wsStatus.Range("E5").Value = 0.8100000000000001
$wsStatus.Range("F5").Value = 1.13
$wsStatus.Range("G5").Value = 47.3

# ---------------------------------------------------------------------------
# Tab "All Trades": columns ... K=Capital After, L=Exit Reason, M=Duration (min)
# Close trade #107 (row 108) and append new trade #140 (row 141)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("G108").Value = 0.01
$wsAll.Range("H108").Value = "CLOSED"
$wsAll.Range("I108").Value = -50
$wsAll.Range("J108").Value = -0.01
$wsAll.Range("K108").Value = 101.13
$wsAll.Range("L108").Value = "early_exit"
$wsAll.Range("M108").Value = 0.13

$wsAll.Range("A141").Value = 140
$wsAll.Range("B141").NumberFormat = "@"
$wsAll.Range("B141").Value = "2026-02-17"
$wsAll.Range("C141").Value = "21:13:41"
$wsAll.Range("D141").Value = "MarketMaking"
$wsAll.Range("E141").Value = "DOWN"
$wsAll.Range("F141").Value = 0.02
$wsAll.Range("G141").Value = ""
$wsAll.Range("H141").Value = "OPEN"
$wsAll.Range("I141").Value = 0
$wsAll.Range("J141").Value = 0
$wsAll.Range("K141").Value = 101.1396151053151
$wsAll.Range("L141").Value = ""
$wsAll.Range("M141").Value = 0
$wsAll.Range("N141").Value = 0
$wsAll.Range("O141").Value = 0
$wsAll.Range("P141").Value = 0.6
$wsAll.Range("Q141").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Tab "MarketMaking": columns ... K=Capital After, L=Entry Slippage, M=Exit
# Slippage, N=Confidence, O=Entry Reason, P=Exit Reason, Q=Duration (min)
# Close trade #107 (row 75) and append new trade #140 (row 108)
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Range("G75").Value = 0.01
$wsMM.Range("H75").Value = "CLOSED"
$wsMM.Range("I75").Value = -50
$wsMM.Range("J75").Value = -0.01
$wsMM.Range("K75").Value = 101.13
$wsMM.Range("P75").Value = "early_exit"
$wsMM.Range("Q75").Value = 0.13

$wsMM.Range("A108").Value = 140
$wsMM.Range("B108").NumberFormat = "@"
$wsMM.Range("B108").Value = "2026-02-17"
$wsMM.Range("C108").Value = "21:13:41"
$wsMM.Range("D108").Value = "MarketMaking"
$wsMM.Range("E108").Value = "DOWN"
$wsMM.Range("F108").Value = 0.02
$wsMM.Range("G108").Value = ""
$wsMM.Range("H108").Value = "OPEN"
$wsMM.Range("I108").Value = 0
$wsMM.Range("J108").Value = 0
$wsMM.Range("K108").Value = 101.1396151053151
$wsMM.Range("L108").Value = 0
$wsMM.Range("M108").Value = 0
$wsMM.Range("N108").Value = 0.6
$wsMM.Range("O108").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("P108").Value = ""
$wsMM.Range("Q108").Value = 0
